$d = $word.ActiveDocument

# Constants
$wdReplaceAll   = 2
$wdFindStop     = 0

# Re-derives the paragraph's Range fresh every call (a Find/Replace collapses
# the range to the replaced text, so reusing a stale Range object would let a
# later search spill into the next paragraph).
function Replace-InParagraph($paraIndex, $findText, $replaceText) {
    $r = $d.Paragraphs($paraIndex).Range
    $r.Find.Execute($findText, $true, $true, $false, $false, $false, $true, $wdFindStop, $false, $replaceText, $wdReplaceAll) | Out-Null
}

# --- "Promoteur" paragraph ---
Replace-InParagraph 34 "Nom, raison sociale, sigle :" "Nom, raison sociale, sigle :Protocole P1 + résumé + 14.1"
Replace-InParagraph 34 "Adresse complète :" "Adresse complète :Protocole P1 + résumé"
Replace-InParagraph 34 "Téléphone :" "Téléphone :Protocole P1 + résumé"

# --- "Investigateur coordonnateur" paragraph ---
Replace-InParagraph 35 "Nom, Prénom :" "Nom, Prénom : HPS I1 ANSM + I1 CPP ; Médicaments G1.3 ANSM; PB G1.1 ANSM HPS I1 ANSM + I1 CPP ; Médicaments G1.1 ANSM; PB G1.3 ANSM"
Replace-InParagraph 35 "Qualité :" "Qualité : HPS I1 ANSM + I1 CPP ; Médicaments G1.4 ANSM; PB G1.4 ANSM"
Replace-InParagraph 35 "Adresse complète :" "Adresse complète : HPS I1 ANSM + I1 CPP"
Replace-InParagraph 35 "Téléphone :" "Téléphone : HPS I1 CPP"
Replace-InParagraph 35 "Mél :" "Mél : HPS I1 ANSM + I1 CPP"

# --- "Recherche" paragraph ---
Replace-InParagraph 37 "Intitulé de la recherche :" "Intitulé de la recherche : Protocole P1 + résumé; DM A ANSM + Q2 CPP; HPS A ANSM + A et Q2 CPP; Médicaments A3 ANSM + Q2 CPP; PB A3 ANSM "
Replace-InParagraph 37 "Numéro d’enregistrement :" "Numéro d’enregistrement : DM A ANSM + Q1 CPP ; HPS A ANSM + Q1 CPP ; PB A2 ANSM"
Replace-InParagraph 37 "Nombre de personnes susceptibles d’être incluses dans la recherche :" "Nombre de personnes susceptibles d’être incluses dans la recherche : Protocole 11.1; Q12 CPP (tous)"

# Insert an extra line break before "Caractéristiques de la recherche :"
# (a new blank line), then refresh the text of that label itself.
$rFind = $d.Paragraphs(37).Range
$found = $rFind.Find.Execute("Caractéristiques de la recherche :", $true, $true, $false, $false, $false, $true, $wdFindStop, $false)
if ($found) {
    $rFind.Collapse(1)
    $rFind.InsertBefore([char]11)
}

Replace-InParagraph 37 "Caractéristiques de la recherche :" "Caractéristiques de la recherche : "
